$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "two drugs can be interacted",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "specified two drugs are interacted",
    2
)
